$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 89.08834838867188
$ws.Range("C2").Value = 7.1724138259887695
$ws.Range("D2").Value = 47.98871994018555
$ws.Range("E2").Value = 57.85714340209961
